$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E14").Value = 45394
$ws.Range("E14").NumberFormat = "DD/MM/YYYY;@"

$ws.Range("A15").Value = $false
$ws.Range("B15").Value = 12
$ws.Range("C13").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D13").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = 45394
$ws.Range("F15").Value = "переместить MainMenu и Boot в project context"

$ws.Range("A16").Value = $false
$ws.Range("B16").Value = 13
$ws.Range("C13").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("D13").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = 45394
$ws.Range("F16").Value = "реализовать Event Bus"

$ws.Range("F17").Select()
